$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.569019626703113
$ws.Range("C2").Value = 4.771454877672043
$ws.Range("D2").Value = 9.111797131699333
$ws.Range("E2").Value = 13.73263851607157
$ws.Range("F2").Value = 33.53764191898374
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 23.12770959703984
$ws.Range("J2").Value = 10.16022615965031
$ws.Range("K2").Value = 9.879565221039517
$ws.Range("M2").Value = 15.2745590891047
$ws.Range("N2").Value = 19.72310094328487
$ws.Range("O2").Value = 25.33748623929653

$ws.Range("B3").Value = 9.323129051461192
$ws.Range("C3").Value = 4.60894037486724
$ws.Range("D3").Value = 9.081071464213856
$ws.Range("E3").Value = 13.72823215512727
$ws.Range("F3").Value = 33.59627219694907
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 23.20840598797368
$ws.Range("J3").Value = 10.1797088264588
$ws.Range("K3").Value = 9.72431153211914
$ws.Range("M3").Value = 15.21705972128948
$ws.Range("N3").Value = 19.78037058264175
$ws.Range("O3").Value = 25.40848426633801

$ws.Range("B4").Value = 9.170404399642287
$ws.Range("C4").Value = 4.50684523000601
$ws.Range("D4").Value = 9.063625518905637
$ws.Range("E4").Value = 13.72782019895621
$ws.Range("F4").Value = 33.639184733143
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 23.26173951049179
$ws.Range("J4").Value = 10.19279432961303
$ws.Range("K4").Value = 9.629292081229258
$ws.Range("M4").Value = 15.18400478770084
$ws.Range("N4").Value = 19.81719153603215
$ws.Range("O4").Value = 25.45667817026678

$ws.Range("B5").Value = 9.10782374012765
$ws.Range("C5").Value = 4.464724283593917
$ws.Range("D5").Value = 9.056878427614548
$ws.Range("E5").Value = 13.72823079897841
$ws.Range("F5").Value = 33.65840857046825
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 23.28442515574758
$ws.Range("J5").Value = 10.19840949986023
$ws.Range("K5").Value = 9.590695588802095
$ws.Range("M5").Value = 15.17111028186017
$ws.Range("N5").Value = 19.83261428344682
$ws.Range("O5").Value = 25.47747334745691

$ws.Range("B6").Value = 9.097414356294637
$ws.Range("C6").Value = 4.457700947733227
$ws.Range("D6").Value = 9.055780105577092
$ws.Range("E6").Value = 13.7283339630114
$ws.Range("F6").Value = 33.66170550223912
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 23.2882495730432
$ws.Range("J6").Value = 10.19935897946641
$ws.Range("K6").Value = 9.584295596422269
$ws.Range("M6").Value = 15.16900420395812
$ws.Range("N6").Value = 19.83520049503986
$ws.Range("O6").Value = 25.48099614154419

$ws.Range("B7").Value = 9.169561679193459
$ws.Range("C7").Value = 4.506279173232699
$ws.Range("D7").Value = 9.063533051743201
$ws.Range("E7").Value = 13.72782339211624
$ws.Range("F7").Value = 33.63943696364402
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 23.26204160338334
$ws.Range("J7").Value = 10.1928689126167
$ws.Range("K7").Value = 9.628770987526702
$ws.Range("M7").Value = 15.18382854448853
$ws.Range("N7").Value = 19.81739783890102
$ws.Range("O7").Value = 25.45695394320454

$ws.Range("B8").Value = 9.484654704773535
$ws.Range("C8").Value = 4.71593647741273
$ws.Range("D8").Value = 9.100911520049603
$ws.Range("E8").Value = 13.73064429063238
$ws.Range("F8").Value = 33.55642183003363
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 23.15474793697541
$ws.Range("J8").Value = 10.1667108830027
$ws.Range("K8").Value = 9.825998445407226
$ws.Range("M8").Value = 15.25427244606801
$ws.Range("N8").Value = 19.74250429799495
$ws.Range("O8").Value = 25.36101079040525

$ws.Range("B9").Value = 10.08476391724618
$ws.Range("C9").Value = 5.106187647836443
$ws.Range("D9").Value = 9.185229100752784
$ws.Range("E9").Value = 13.75428849249345
$ws.Range("F9").Value = 33.44854839740488
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 22.97438779415045
$ws.Range("J9").Value = 10.12431202030644
$ws.Range("K9").Value = 10.21313256478055
$ws.Range("M9").Value = 15.4098052290229
$ws.Range("N9").Value = 19.60873130666219
$ws.Range("O9").Value = 25.20942116238809

$ws.Range("B10").Value = 10.50994793123327
$ws.Range("C10").Value = 5.377148936399003
$ws.Range("D10").Value = 9.253543277860075
$ws.Range("E10").Value = 13.78257224336038
$ws.Range("F10").Value = 33.40284094134071
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 22.86019884663679
$ws.Range("J10").Value = 10.09856623543422
$ws.Range("K10").Value = 10.4950804859546
$ws.Range("M10").Value = 15.53402794524838
$ws.Range("N10").Value = 19.5183506411426
$ws.Range("O10").Value = 25.1203967641841

$ws.Range("B11").Value = 10.69904819112866
$ws.Range("C11").Value = 5.496463380679355
$ws.Range("D11").Value = 9.285922184016757
$ws.Range("E11").Value = 13.79777566831
$ws.Range("F11").Value = 33.38933695969613
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 22.8122314990136
$ws.Range("J11").Value = 10.08802319577594
$ws.Range("K11").Value = 10.62225782871483
$ws.Range("M11").Value = 15.59255218842386
$ws.Range("N11").Value = 19.47893309572037
$ws.Range("O11").Value = 25.0847620762475

$ws.Range("B12").Value = 10.76996648470059
$ws.Range("C12").Value = 5.541039428222196
$ws.Range("D12").Value = 9.298363194873556
$ws.Range("E12").Value = 13.80386558082744
$ws.Range("F12").Value = 33.38527095118867
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 22.79463976348419
$ws.Range("J12").Value = 10.08419855131265
$ws.Range("K12").Value = 10.67021719755998
$ws.Range("M12").Value = 15.61499019296583
$ws.Range("N12").Value = 19.46424948458602
$ws.Range("O12").Value = 25.07196828513259

$ws.Range("B13").Value = 10.75472467813306
$ws.Range("C13").Value = 5.531466678384333
$ws.Range("D13").Value = 9.295675924594846
$ws.Range("E13").Value = 13.80253926835032
$ws.Range("F13").Value = 33.38610005646007
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 22.79840299089055
$ws.Range("J13").Value = 10.085014799973
$ws.Range("K13").Value = 10.65989786709581
$ws.Range("M13").Value = 15.61014570590881
$ws.Range("N13").Value = 19.46740107442283
$ws.Range("O13").Value = 25.07469249643872

$ws.Range("B14").Value = 10.70489688164824
$ws.Range("C14").Value = 5.500143021202121
$ws.Range("D14").Value = 9.286942153257513
$ws.Range("E14").Value = 13.79827004134585
$ws.Range("F14").Value = 33.388981454449
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 22.81077274078306
$ws.Range("J14").Value = 10.08770517926294
$ws.Range("K14").Value = 10.62620772897551
$ws.Range("M14").Value = 15.59439272136012
$ws.Range("N14").Value = 19.47772020310201
$ws.Range("O14").Value = 25.0836954830889

$ws.Range("B15").Value = 10.67428413135373
$ws.Range("C15").Value = 5.480876437148758
$ws.Range("D15").Value = 9.281615658202961
$ws.Range("E15").Value = 13.79569824040551
$ws.Range("F15").Value = 33.390882808054
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 22.81842413441557
$ws.Range("J15").Value = 10.08937495448724
$ws.Range("K15").Value = 10.60554422617436
$ws.Range("M15").Value = 15.58477911602487
$ws.Range("N15").Value = 19.48407257581257
$ws.Range("O15").Value = 25.0893012990236

$ws.Range("B16").Value = 10.49749686503972
$ws.Range("C16").Value = 5.369268663383861
$ws.Range("D16").Value = 9.251452816277526
$ws.Range("E16").Value = 13.78162540996283
$ws.Range("F16").Value = 33.40387009788871
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 22.8634137076868
$ws.Range("J16").Value = 10.09927874236584
$ws.Range("K16").Value = 10.48674354130715
$ws.Range("M16").Value = 15.53024263877538
$ws.Range("N16").Value = 19.52096072738759
$ws.Range("O16").Value = 25.12282351827059

$ws.Range("B17").Value = 10.38788669905488
$ws.Range("C17").Value = 5.29976132229846
$ws.Range("D17").Value = 9.23327724733762
$ws.Range("E17").Value = 13.77358844924831
$ws.Range("F17").Value = 33.41370401388642
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 22.89203248174945
$ws.Range("J17").Value = 10.10565355243559
$ws.Range("K17").Value = 10.41355389128371
$ws.Range("M17").Value = 15.49729276797787
$ws.Range("N17").Value = 19.54402430340318
$ws.Range("O17").Value = 25.14463458227147

$ws.Range("B18").Value = 10.32443913247778
$ws.Range("C18").Value = 5.259412982566822
$ws.Range("D18").Value = 9.222946197332552
$ws.Range("E18").Value = 13.76918599020669
$ws.Range("F18").Value = 33.42004631514898
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 22.90886762295332
$ws.Range("J18").Value = 10.10943021065524
$ws.Range("K18").Value = 10.37135782060433
$ws.Range("M18").Value = 15.4785313467413
$ws.Range("N18").Value = 19.55744966264968
$ws.Range("O18").Value = 25.15763740044087

$ws.Range("B19").Value = 10.30288991087331
$ws.Range("C19").Value = 5.245689494512289
$ws.Range("D19").Value = 9.219469639627162
$ws.Range("E19").Value = 13.76773330839859
$ws.Range("F19").Value = 33.42231155119774
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 22.91463199866654
$ws.Range("J19").Value = 10.11072783012587
$ws.Range("K19").Value = 10.357055277956
$ws.Range("M19").Value = 15.47221217250771
$ws.Range("N19").Value = 19.56202273753592
$ws.Range("O19").Value = 25.16211850437138

$ws.Range("B20").Value = 10.39959708430376
$ws.Range("C20").Value = 5.307199035877354
$ws.Range("D20").Value = 9.235199385066815
$ws.Range("E20").Value = 13.77442123220477
$ws.Range("F20").Value = 33.41258617017787
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 22.88894721423664
$ws.Range("J20").Value = 10.10496355686757
$ws.Range("K20").Value = 10.42135566053333
$ws.Range("M20").Value = 15.50078072173956
$ws.Range("N20").Value = 19.54155261547343
$ws.Range("O20").Value = 25.14226538241694

$ws.Range("B21").Value = 10.7195517562101
$ws.Range("C21").Value = 5.509360259677662
$ws.Range("D21").Value = 9.289502654084956
$ws.Range("E21").Value = 13.79951501535752
$ws.Range("F21").Value = 33.38810669092835
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 22.80712390535857
$ws.Range("J21").Value = 10.08691039871988
$ws.Range("K21").Value = 10.63610910919996
$ws.Range("M21").Value = 15.59901237042117
$ws.Range("N21").Value = 19.47468263961102
$ws.Range("O21").Value = 25.08103207447999

$ws.Range("B22").Value = 10.92461170867835
$ws.Range("C22").Value = 5.63793901995797
$ws.Range("D22").Value = 9.326037479057279
$ws.Range("E22").Value = 13.81785276707409
$ws.Range("F22").Value = 33.37821391088993
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 22.75698458370277
$ws.Range("J22").Value = 10.07608939538155
$ws.Range("K22").Value = 10.77527834090967
$ws.Range("M22").Value = 15.66481544439945
$ws.Range("N22").Value = 19.4323948957522
$ws.Range("O22").Value = 25.04509479418706

$ws.Range("B23").Value = 10.81555907923549
$ws.Range("C23").Value = 5.569649949150634
$ws.Range("D23").Value = 9.306445115785985
$ws.Range("E23").Value = 13.80788945962646
$ws.Range("F23").Value = 33.38293545343486
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 22.7834393997558
$ws.Range("J23").Value = 10.08177540404343
$ws.Range("K23").Value = 10.70112364271037
$ws.Range("M23").Value = 15.62955301259028
$ws.Range("N23").Value = 19.4548354821587
$ws.Range("O23").Value = 25.06390136906759

$ws.Range("B24").Value = 10.39430415698705
$ws.Range("C24").Value = 5.303837650045877
$ws.Range("D24").Value = 9.234330017470459
$ws.Range("E24").Value = 13.77404405144015
$ws.Range("F24").Value = 33.41308940219238
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 22.89034087434274
$ws.Range("J24").Value = 10.10527515592272
$ws.Range("K24").Value = 10.4178288459405
$ws.Range("M24").Value = 15.49920325022397
$ws.Range("N24").Value = 19.54266954903261
$ws.Range("O24").Value = 25.14333505444738

$ws.Range("B25").Value = 9.924848208738471
$ws.Range("C25").Value = 5.003183817868677
$ws.Range("D25").Value = 9.16127430397159
$ws.Range("E25").Value = 13.7459653219052
$ws.Range("F25").Value = 33.47184279513991
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 23.01996306458391
$ws.Range("J25").Value = 10.12431202030644
$ws.Range("K25").Value = 10.21313256478055
$ws.Range("M25").Value = 15.36593250398711
$ws.Range("N25").Value = 19.64352700501354
$ws.Range("O25").Value = 25.20942116238809
